$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 1990
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 1990
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 1990
$ws.Range("M13").Value = $null
$ws.Range("N13").Value = -2328

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1991.6666
$ws.Range("I98").Value = 485.2
$ws.Range("J98").Value = 3874.75
$ws.Range("K98").Value = 485.2
$ws.Range("L98").Value = 3874.75
$ws.Range("M98").Value = 1012.8
$ws.Range("N98").Value = -6870.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 16667299
$ws.Range("I101").Value = 25000548
$ws.Range("J101").Value = 800
$ws.Range("K101").Value = 75001644
$ws.Range("L101").Value = 2400
$ws.Range("M101").Value = -75000022
$ws.Range("N101").Value = -5644

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H115").Value = 185
$ws.Range("I115").Value = 185
$ws.Range("K115").Value = 555
$ws.Range("M115").Value = 1012

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1991.6666
$ws.Range("I122").Value = 485.2
$ws.Range("J122").Value = 3874.75
$ws.Range("K122").Value = 1455.6
$ws.Range("L122").Value = 11624.25
$ws.Range("M122").Value = 994.4000000000001
$ws.Range("N122").Value = -16524.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2332.3333
$ws.Range("I137").Value = 1000
$ws.Range("J137").Value = 2998.5
$ws.Range("K137").Value = 3000
$ws.Range("L137").Value = 8995.5
$ws.Range("M137").Value = -450
$ws.Range("N137").Value = -14095.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6418.067
$ws.Range("I61").Value = 4129.8335
$ws.Range("J61").Value = 7943.5557
$ws.Range("K61").Value = 4129.8335
$ws.Range("L61").Value = 7943.5557
$ws.Range("M61").Value = -3917.8335
$ws.Range("N61").Value = -8367.555700000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1662.5454
$ws.Range("I132").Value = 1662.5454
$ws.Range("K132").Value = 4987.6362
$ws.Range("M132").Value = -2457.6362

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 6418.067
$ws.Range("I136").Value = 4129.8335
$ws.Range("J136").Value = 7943.5557
$ws.Range("K136").Value = 12389.5005
$ws.Range("L136").Value = 23830.6671
$ws.Range("M136").Value = -9839.500499999998
$ws.Range("N136").Value = -28930.6671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 788.44446
$ws.Range("I20").Value = 807.8333
$ws.Range("J20").Value = 749.6667
$ws.Range("K20").Value = 807.8333
$ws.Range("L20").Value = 749.6667
$ws.Range("M20").Value = -560.8333
$ws.Range("N20").Value = -1243.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1288.4445
$ws.Range("I94").Value = 1288.4445
$ws.Range("K94").Value = 1288.4445
$ws.Range("M94").Value = -837.4445000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2081.4285
$ws.Range("I134").Value = 2081.4285
$ws.Range("K134").Value = 6244.2855
$ws.Range("M134").Value = -3709.2855

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1655.5
$ws.Range("I31").Value = 1655.5
$ws.Range("K31").Value = 1655.5
$ws.Range("M31").Value = -1360.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1655.5
$ws.Range("I34").Value = 1655.5
$ws.Range("K34").Value = 1655.5
$ws.Range("M34").Value = -1453.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1794.4166
$ws.Range("I58").Value = 1928.8889
$ws.Range("J58").Value = 1391
$ws.Range("K58").Value = 1928.8889
$ws.Range("L58").Value = 1391
$ws.Range("M58").Value = -1725.8889
$ws.Range("N58").Value = -1797

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1447.2609
$ws.Range("I132").Value = 1515.1052
$ws.Range("J132").Value = 1125
$ws.Range("K132").Value = 4545.3156
$ws.Range("L132").Value = 3375
$ws.Range("M132").Value = -2015.3156
$ws.Range("N132").Value = -8435

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 966.2308
$ws.Range("I134").Value = 903.9167
$ws.Range("J134").Value = 1714
$ws.Range("K134").Value = 2711.7501
$ws.Range("L134").Value = 5142
$ws.Range("M134").Value = -176.7501000000002
$ws.Range("N134").Value = -10212

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1794.4166
$ws.Range("I136").Value = 1928.8889
$ws.Range("J136").Value = 1391
$ws.Range("K136").Value = 5786.6667
$ws.Range("L136").Value = 4173
$ws.Range("M136").Value = -3236.6667
$ws.Range("N136").Value = -9273

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 59978
$ws.Range("J37").Value = 59978
$ws.Range("L37").Value = 179934
$ws.Range("N37").Value = -180158

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 2950
$ws.Range("I104").Value = 2950
$ws.Range("K104").Value = 8850
$ws.Range("M104").Value = -6229

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 9688.6875
$ws.Range("I140").Value = 3337.3333
$ws.Range("J140").Value = 13499.5
$ws.Range("K140").Value = 10011.9999
$ws.Range("L140").Value = 40498.5
$ws.Range("M140").Value = -4831.999899999999
$ws.Range("N140").Value = -50858.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 7439.4
$ws.Range("I141").Value = 7439.4
$ws.Range("K141").Value = 22318.2
$ws.Range("M141").Value = -17138.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2815.8572
$ws.Range("I126").Value = 2906
$ws.Range("J126").Value = 2779.8
$ws.Range("K126").Value = 8718
$ws.Range("L126").Value = 8339.400000000001
$ws.Range("M126").Value = -6248
$ws.Range("N126").Value = -13279.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1736.6666
$ws.Range("I61").Value = 816.6667
$ws.Range("J61").Value = 2656.6667
$ws.Range("K61").Value = 816.6667
$ws.Range("L61").Value = 2656.6667
$ws.Range("M61").Value = -614.6667
$ws.Range("N61").Value = -3060.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1736.6666
$ws.Range("I113").Value = 816.6667
$ws.Range("J113").Value = 2656.6667
$ws.Range("K113").Value = 816.6667
$ws.Range("L113").Value = 2656.6667
$ws.Range("M113").Value = 1353.3333
$ws.Range("N113").Value = -6996.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4183.8096
$ws.Range("I136").Value = 3939.8948
$ws.Range("J136").Value = 6501
$ws.Range("K136").Value = 11819.6844
$ws.Range("L136").Value = 19503
$ws.Range("M136").Value = -9269.6844
$ws.Range("N136").Value = -24603

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 11967.125
$ws.Range("I62").Value = 10958.5
$ws.Range("J62").Value = 14993
$ws.Range("K62").Value = 10958.5
$ws.Range("L62").Value = 14993
$ws.Range("M62").Value = -10334.5
$ws.Range("N62").Value = -16241

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 11967.125
$ws.Range("I65").Value = 10958.5
$ws.Range("J65").Value = 14993
$ws.Range("K65").Value = 54792.5
$ws.Range("L65").Value = 74965
$ws.Range("M65").Value = -51672.5
$ws.Range("N65").Value = -81205

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 938.82355
$ws.Range("I107").Value = 1072.6666
$ws.Range("J107").Value = 788.25
$ws.Range("K107").Value = 3217.9998
$ws.Range("L107").Value = 2364.75
$ws.Range("M107").Value = -1297.9998
$ws.Range("N107").Value = -6204.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 604.5714
$ws.Range("I113").Value = 538.6667
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 1616.0001
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 553.9999
$ws.Range("N113").Value = -7340

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2256.3948
$ws.Range("I132").Value = 2370.4688
$ws.Range("J132").Value = 1648
$ws.Range("K132").Value = 7111.4064
$ws.Range("L132").Value = 4944
$ws.Range("M132").Value = -4581.4064
$ws.Range("N132").Value = -10004

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3189.913
$ws.Range("I136").Value = 3384
$ws.Range("J136").Value = 1896
$ws.Range("K136").Value = 10152
$ws.Range("L136").Value = 5688
$ws.Range("M136").Value = -7602
$ws.Range("N136").Value = -10788
